# Generate Report for Handback
# A new handback/handoff cycle was recorded for the "10c79ea7..." source
# file (row 2 of each status sheet): its latest HO xliff generation date,
# correspond-handoff datetime and correspond-handback datetime all moved
# forward. The "84af6fea..." row (row 3) is untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-17 08:47:28"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-17 08:47:23"
$zhcn.Range("K2").Value = "2016-08-17 08:47:40"

# --- de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-17 08:47:28"
$dede.Range("K2").Value = "2016-08-17 08:47:47"
